# Trading update: 2026-02-17 20:49:49
# Appends the latest MarketMaking trade (Trade #74) to the "All Trades"
# log (row 75) and to the per-strategy "MarketMaking" sheet (row 42).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "All Trades" sheet -> new row 75
#   Trade #, Date, Time, Strategy, Side, Entry, Exit, Status, P&L%, P&L$,
#   Capital After, Exit Reason, Duration, Entry Slip, Exit Slip,
#   Confidence, Entry Reason
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$rowAT = 75

$allTrades.Cells.Item($rowAT, 1).Value = 74

# Date column would otherwise be auto-parsed into a date serial number;
# force text storage, write it, then drop back to the default ("Normal")
# style so no extra formatting sticks to the cell.
$allTrades.Cells.Item($rowAT, 2).NumberFormat = "@"
$allTrades.Cells.Item($rowAT, 2).Value = "2026-02-17"
$allTrades.Cells.Item($rowAT, 2).Style = "Normal"

$allTrades.Cells.Item($rowAT, 3).Value = "20:49:20"
$allTrades.Cells.Item($rowAT, 4).Value = "MarketMaking"
$allTrades.Cells.Item($rowAT, 5).Value = "UP"
$allTrades.Cells.Item($rowAT, 6).Value = 0.03
# G75 (Exit Price) stays blank - trade is still OPEN
$allTrades.Cells.Item($rowAT, 8).Value = "OPEN"
$allTrades.Cells.Item($rowAT, 9).Value = 0
$allTrades.Cells.Item($rowAT, 10).Value = 0
$allTrades.Cells.Item($rowAT, 11).Value = 100.305976116214
# L75 (Exit Reason) stays blank - trade is still OPEN
$allTrades.Cells.Item($rowAT, 13).Value = 0
$allTrades.Cells.Item($rowAT, 14).Value = 0
$allTrades.Cells.Item($rowAT, 15).Value = 0
$allTrades.Cells.Item($rowAT, 16).Value = 0.6
$allTrades.Cells.Item($rowAT, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# "MarketMaking" sheet -> new row 42
#   Trade #, Date, Time, Strategy, Side, Entry, Exit, Status, P&L%, P&L$,
#   Capital After, Entry Slip, Exit Slip, Confidence, Entry Reason,
#   Exit Reason, Duration
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$rowMM = 42

$mm.Cells.Item($rowMM, 1).Value = 74

$mm.Cells.Item($rowMM, 2).NumberFormat = "@"
$mm.Cells.Item($rowMM, 2).Value = "2026-02-17"
$mm.Cells.Item($rowMM, 2).Style = "Normal"

$mm.Cells.Item($rowMM, 3).Value = "20:49:20"
$mm.Cells.Item($rowMM, 4).Value = "MarketMaking"
$mm.Cells.Item($rowMM, 5).Value = "UP"
$mm.Cells.Item($rowMM, 6).Value = 0.03
# G42 (Exit Price) stays blank - trade is still OPEN
$mm.Cells.Item($rowMM, 8).Value = "OPEN"
$mm.Cells.Item($rowMM, 9).Value = 0
$mm.Cells.Item($rowMM, 10).Value = 0
$mm.Cells.Item($rowMM, 11).Value = 100.305976116214
$mm.Cells.Item($rowMM, 12).Value = 0
$mm.Cells.Item($rowMM, 13).Value = 0
$mm.Cells.Item($rowMM, 14).Value = 0.6
$mm.Cells.Item($rowMM, 15).Value = "Normal spread capture: 19600 bps"
# P42 (Exit Reason) stays blank - trade is still OPEN
$mm.Cells.Item($rowMM, 17).Value = 0
